# Fruta / hortaliza, semanal
# The weekly refresh rotates the "Calidad/Volumen/Precio/..." data among
# three groups of rows (by date) while keeping the fixed descriptive
# columns (A,B,C,E,F,G,H,I,J,K) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that get rotated between rows: D,L,M,N,O,P,Q,R,S,T
$cols = @("D","L","M","N","O","P","Q","R","S","T")

function Get-RowValues($row) {
    $vals = @{}
    foreach ($c in $cols) {
        $vals[$c] = $ws.Range("$c$row").Value2
    }
    return $vals
}

function Set-RowValues($row, $vals) {
    foreach ($c in $cols) {
        $ws.Range("$c$row").Value2 = $vals[$c]
    }
}

# Capture original values before any writes (since rows feed each other).
$orig2  = Get-RowValues 2
$orig3  = Get-RowValues 3
$orig8  = Get-RowValues 8
$orig9  = Get-RowValues 9
$orig10 = Get-RowValues 10
$orig11 = Get-RowValues 11
$orig12 = Get-RowValues 12
$orig13 = Get-RowValues 13

# Cycle (2 -> 11 -> 9 -> 2): after[2]=before[11], after[11]=before[9], after[9]=before[2]
Set-RowValues 2  $orig11
Set-RowValues 11 $orig9
Set-RowValues 9  $orig2

# Cycle (3 -> 12 -> 3): after[3]=before[12], after[12]=before[3]
Set-RowValues 3  $orig12
Set-RowValues 12 $orig3

# Cycle (8 -> 13 -> 10 -> 8): after[8]=before[13], after[13]=before[10], after[10]=before[8]
Set-RowValues 8  $orig13
Set-RowValues 13 $orig10
Set-RowValues 10 $orig8
